# Auto-generated Excel COM-interop script applying the "Update countries & provincias Spain" commit.
# Updates country-name relabeling (shared-string reorder in the source) and refreshed COVID case counters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name relabels (rows whose displayed country changed) ---
$ws.Range("A51").Value = "Portugal"
$ws.Range("A52").Value = "Nepal"
$ws.Range("A71").Value = "Kenia"
$ws.Range("A72").Value = "Paraguay"
$ws.Range("A96").Value = "Namibia"
$ws.Range("A97").Value = "Malasia"
$ws.Range("A98").Value = "Birmania"
$ws.Range("A99").Value = "Consejo Danes para los Refugiados"
$ws.Range("A100").Value = "Guinea"
$ws.Range("A101").Value = "Montenegro"
$ws.Range("A102").Value = "Maldivas"
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Updated statistic counters (Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes) ---
$ws.Range("B4").Value = 7289426
$ws.Range("C4").Value = 1865
$ws.Range("D4").Value = 4524425
$ws.Range("E4").Value = 2555795
$ws.Range("G4").Value = 29
$ws.Range("H4").Value = 209206
$ws.Range("B5").Value = 6005795
$ws.Range("C5").Value = 15214
$ws.Range("D5").Value = 4945998
$ws.Range("E5").Value = 965215
$ws.Range("G5").Value = 48
$ws.Range("H5").Value = 94582
$ws.Range("B20").Value = 333193
$ws.Range("C20").Value = 403
$ws.Range("D20").Value = 317005
$ws.Range("E20").Value = 11505
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 4683
$ws.Range("B25").Value = 285270
$ws.Range("C25").Value = 245
$ws.Range("E25").Value = 26238
$ws.Range("B51").Value = 73604
$ws.Range("C51").Value = 665
$ws.Range("D51").Value = 47647
$ws.Range("E51").Value = 24004
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 1953
$ws.Range("B52").Value = 73394
$ws.Range("C52").Value = 1573
$ws.Range("D52").Value = 53898
$ws.Range("E52").Value = 19019
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 477
$ws.Range("E56").Value = 6281
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 242
$ws.Range("D59").Value = 57367
$ws.Range("E59").Value = 306
$ws.Range("B60").Value = 55320
$ws.Range("C60").Value = 501
$ws.Range("D60").Value = 51829
$ws.Range("E60").Value = 3033
$ws.Range("G60").Value = 6
$ws.Range("H60").Value = 458
$ws.Range("B70").Value = 38703
$ws.Range("C70").Value = 450
$ws.Range("D70").Value = 29068
$ws.Range("E70").Value = 9344
$ws.Range("G70").Value = 6
$ws.Range("H70").Value = 291
$ws.Range("B71").Value = 38115
$ws.Range("C71").Value = 244
$ws.Range("D71").Value = 24621
$ws.Range("E71").Value = 12803
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 691
$ws.Range("B72").Value = 37922
$ws.Range("D72").Value = 21757
$ws.Range("E72").Value = 15383
$ws.Range("H72").Value = 782
$ws.Range("B75").Value = 33384
$ws.Range("C75").Value = 72
$ws.Range("E75").Value = 1101
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 747
$ws.Range("B76").Value = 32364
$ws.Range("C76").Value = 536
$ws.Range("D76").Value = 18128
$ws.Range("E76").Value = 13716
$ws.Range("G76").Value = 21
$ws.Range("H76").Value = 520
$ws.Range("B93").Value = 13660
$ws.Range("C93").Value = 33
$ws.Range("E93").Value = 2200
$ws.Range("B96").Value = 11033
$ws.Range("C96").Value = 115
$ws.Range("D96").Value = 8776
$ws.Range("E96").Value = 2137
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 120
$ws.Range("B97").Value = 10919
$ws.Range("C97").Value = 150
$ws.Range("D97").Value = 9835
$ws.Range("E97").Value = 950
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 134
$ws.Range("B98").Value = 10734
$ws.Range("C98").Value = 743
$ws.Range("D98").Value = 2862
$ws.Range("E98").Value = 7646
$ws.Range("G98").Value = 28
$ws.Range("H98").Value = 226
$ws.Range("B99").Value = 10593
$ws.Range("D99").Value = 10093
$ws.Range("E99").Value = 229
$ws.Range("H99").Value = 271
$ws.Range("B100").Value = 10512
$ws.Range("D100").Value = 9836
$ws.Range("E100").Value = 611
$ws.Range("H100").Value = 65
$ws.Range("B101").Value = 10197
$ws.Range("D101").Value = 6368
$ws.Range("E101").Value = 3671
$ws.Range("H101").Value = 158
$ws.Range("B102").Value = 10045
$ws.Range("D102").Value = 8754
$ws.Range("E102").Value = 1257
$ws.Range("H102").Value = 34
$ws.Range("B105").Value = 9646
$ws.Range("C105").Value = 41
$ws.Range("D105").Value = 8430
$ws.Range("E105").Value = 1141
$ws.Range("B142").Value = 3352
$ws.Range("C142").Value = 3
$ws.Range("E142").Value = 131
$ws.Range("B149").Value = 2623
$ws.Range("C149").Value = 22
$ws.Range("D149").Value = 2158
$ws.Range("E149").Value = 455

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 16:14"
